$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 16: "zweiwöchentliches Meeting" ---------------------------
# Copy the formatting (number format / wrap text) from the row above so the
# new cells pick up the same style ids instead of inventing new ones.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("D16:D17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A16").Value2 = 44173
$ws.Range("B16").Value = 1
$ws.Range("D16").Value = "zweiwöchentliches Meeting"

# Extend the running-total formula from C10 down through the new C16 in one
# shot so it is filled as a single shared formula (same as dragging the
# fill handle down in Excel).
$ws.Range("C10:C16").Formula = "=B10+C9"

# --- New row 17: "redpitaya abgeholt, ..." ------------------------------
$ws.Range("A17").Value2 = 44182
$ws.Range("B17").Value = 4
$ws.Range("C17").Formula = "=C16+B17"
$ws.Range("D17").Value = "redpitaya abgeholt, upgedated und in HomeOffice Arbeitsplatz integriert "

$ws.Rows("17").RowHeight = 30

$ws.Calculate()

# Leave the view roughly where the author would have been scrolled to/
# selected after adding the new entries.
[void]$ws.Activate()
[void]$ws.Range("D27").Select()
